$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Header style: font color -> white, on both sheets (shared style) ---
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# --- Sheet1: Training Dashboard ---

# Column J (STATUS) width 8 -> 11
$ws1.Columns.Item(10).ColumnWidth = 10.16666666666667

# Row 3
$ws1.Range("H3").Value = 394
$ws1.Range("I3").Value = "'16-Sep-2025"

# Row 4
$ws1.Range("H4").Value = 413
$ws1.Range("I4").Value = "'16-Sep-2025"

# Row 5
$ws1.Range("H5").Value = 413
$ws1.Range("I5").Value = "'16-Sep-2025"

# Row 6: highlight pink (NOT VALID) + data updates
$ws1.Range("A6:K6").Interior.Color = 13551615
$ws1.Range("H6").Value = 14
$ws1.Range("I6").Value = "'16-Sep-2025"
$ws1.Range("J6").Value = "NOT VALID"

# --- Sheet2: Exam Dashboard ---

# Column E width 44 -> 15
$ws2.Columns.Item(5).ColumnWidth = 14.16666666666667

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"

Write-Host "edits applied"
